$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row 2 for the "2022-Q4" quarter,
#    pushing the existing quarters down by one row, and renumber the index
#    column (A) so it stays 0,1,2,...
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()

# Rows 3..9 now hold what used to be rows 2..8 - fix up formatting that the
# row-insert copied down from row 1 (header) onto B2:D2, and restore column A
# numbering (0-based) for every data row.
$summary.Range("B2:D2").ClearFormats()

$summary.Cells.Item(3, 1).Copy()
$summary.Cells.Item(2, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q4"
$summary.Cells.Item(2, 3).Value = 6
$summary.Cells.Item(2, 4).Value = 0.36

for ($r = 3; $r -le 9; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------------
# 2) Insert a brand new "2022-Q4" fund-holdings sheet right after "总计"
#    (i.e. before the old "2022-Q2" sheet), matching the layout used by the
#    other quarterly sheets.
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")

$newSheet = $wb.Worksheets.Add($null, $summary)
$newSheet.Name = "2022-Q4"

# Pull header (row 1) and column-A formatting from the existing "2022-Q2"
# sheet so the new sheet matches the established look (bold + border style).
$q2.Activate()
$q2.Range("B1:H1").Copy()
$newSheet.Activate()
$newSheet.Range("B1:H1").Select()
$newSheet.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$q2.Activate()
$q2.Range("A2").Copy()
$newSheet.Activate()
$newSheet.Range("A2:A7").Select()
$newSheet.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Columns B-G hold text (fund code keeps leading zeros, numeric columns are
# stored as plain text strings like the other quarterly sheets).
$newSheet.Range("B2:G7").NumberFormat = "@"

$data = @(
    @(0, "008347", "中信建投价值甄选混合A", "4.72", "63.36", "2.05", "0.0968", 9),
    @(1, "003822", "中信建投行业轮换混合A", "4.37", "56.48", "1.93", "0.0843", 7),
    @(2, "003823", "中信建投行业轮换混合C", "2.72", "56.48", "1.93", "0.0525", 7),
    @(3, "013340", "创金合信芯片产业股票C", "0.94", "92.41", "4.46", "0.0419", 8),
    @(4, "013339", "创金合信芯片产业股票A", "0.92", "92.41", "4.46", "0.0410", 8),
    @(5, "008348", "中信建投价值甄选混合C", "1.98", "63.36", "2.05", "0.0406", 9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rd = $data[$i]
    $newSheet.Cells.Item($row, 1).Value = $rd[0]
    $newSheet.Cells.Item($row, 2).Value = $rd[1]
    $newSheet.Cells.Item($row, 3).Value = $rd[2]
    $newSheet.Cells.Item($row, 4).Value = $rd[3]
    $newSheet.Cells.Item($row, 5).Value = $rd[4]
    $newSheet.Cells.Item($row, 6).Value = $rd[5]
    $newSheet.Cells.Item($row, 7).Value = $rd[6]
    $newSheet.Cells.Item($row, 8).Value = $rd[7]
}

# Drop the temporary "@" number-format style (keep the underlying text type)
# so the data cells end up unstyled, same as the sibling quarterly sheets.
$newSheet.Range("B2:G7").ClearFormats()

$summary.Activate()
$summary.Range("A1").Select()

Write-Output "done"
